# "Out of PO.xlsx" - re-upload of the player table with rows reordered.
# The header row (A1:C1) and the overall table shape (A1:C18) are unchanged;
# only the order of the 17 player rows (A2:C18) is different in the new file.
# Same 17 players / positions / teams as before, just reshuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Payton Pritchard",
    "Bradley Beal",
    "Duncan Robinson",
    "Dyson Daniels",
    "Michael Porter Jr.",
    "De'Andre Hunter",
    "Kristaps Porzingis",
    "Alperen Sengün",
    "Deandre Ayton",
    "Malik Beasley",
    "Ziaire Williams",
    "Victor Wembanyama",
    "Josh Hart",
    "Domantas Sabonis",
    "Donovan Mitchell",
    "Tari Eason",
    "Cam Thomas"
)

$positions = @(
    "PG,SG",
    "PG,SG,SF",
    "SG,SF",
    "PG,SG,SF",
    "SF,PF",
    "SF,PF",
    "PF,C",
    "C",
    "C",
    "SG,SF",
    "SG,SF",
    "C",
    "SG,SF,PF",
    "C",
    "PG,SG",
    "SF,PF",
    "SG,SF"
)

$teams = @(
    "Boston Celtics",
    "Phoenix Suns",
    "Miami Heat",
    "Atlanta Hawks",
    "Denver Nuggets",
    "Atlanta Hawks",
    "Boston Celtics",
    "Houston Rockets",
    "Portland Trail Blazers",
    "Detroit Pistons",
    "Brooklyn Nets",
    "San Antonio Spurs",
    "New York Knicks",
    "Sacramento Kings",
    "Cleveland Cavaliers",
    "Houston Rockets",
    "Brooklyn Nets"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
